$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AIC_PinsFrames1")
$ws.Cells.Item(2, 4).Value = "pass"
